# Recreate the author's edit: add a SUM(A:A) formula in C2, switch the
# workbook to automatic calculation (dropping the stale "manual" calc
# mode left over from earlier testing), and leave the selection on the
# newly entered cell, just like Excel does after you type a formula and
# press Enter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook was left in manual calculation mode; switch back to
# automatic so the new formula (and everything else) recalculates normally.
$excel.Calculation = -4105   # xlCalculationAutomatic

# Enter the new formula in C2 - Excel will calculate its value (21).
$ws.Range("C2").Formula = "=SUM(A:A)"

# Match the author's final selection state (C2 selected/active).
$ws.Range("C2").Select()
